$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "in-dev" -> "open", "completed" -> "closed"
$ws.Range("G1").Value = "open"
$ws.Range("H1").Value = "closed"

# Update the active selection to G2 (as seen in sheetView selection)
$ws.Range("G2").Select()
